$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in contracted_rte (column E) for rows 3 through 11 with 0.832
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = 0.832
}

# Add new row 14: 2025-05
$ws.Cells.Item(14, 1).Value = "2025-05"
$ws.Cells.Item(14, 2).Value = 0.8023110536204735
$ws.Cells.Item(14, 3).Value = 1007
$ws.Cells.Item(14, 5).Value = 0.832

# Add new row 15: 2025-06
$ws.Cells.Item(15, 1).Value = "2025-06"
$ws.Cells.Item(15, 2).Value = 0.8902329256087044
$ws.Cells.Item(15, 3).Value = 97
$ws.Cells.Item(15, 5).Value = 0.832
